# :ambulance: Hotfix for excel fortune
# Update wording of a few fortune messages and restore the intended ordering
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wording fixes (new message text), in authoring order so the
# shared-string table is appended to in the same order as upstream
$ws.Range("A6").Value = "조금만 더 노력한다면 이번엔 좋은 점수 받을 것 같아요. 화이팅"
$ws.Range("A24").Value = "일 년 전 당신을 떠올려 보세요. 분명 계속 성장했어요."
$ws.Range("A22").Value = "내일의 나는 오늘의 나보다 나을 거에요."
$ws.Range("A21").Value = "할 수 있다고 믿어 봐요. 생각보다 별 것 아닐 수도 있어요."
$ws.Range("A20").Value = "당신은 웃는 얼굴이 참 예뻐요. 당신의 미소가 힘이 돼요."

# Re-sequence the trailing block of messages (rows 26-34)
$ws.Range("A26").Value = "무엇보다도 나 자신을 사랑하는 것이 중요해요."
$ws.Range("A27").Value = "행복은 거창한 것이 아니라 작은 것에서 시작해요."
$ws.Range("A28").Value = "기분 전환하러 이번 방학엔 여행을 떠나보는 것은 어떨까요"
$ws.Range("A29").Value = "오지 않은 미래를 걱정하는 것 보다 마주한 현재에 최선을 다하세요"
$ws.Range("A30").Value = "모든 것은 나의 태도에 달려있어요"
$ws.Range("A31").Value = "우리는 아직 우리가 보석인지 몰라요. 당신은 특별해요"
$ws.Range("A32").Value = "익숙함에 속아 소중함을 잃지 말자구요"
$ws.Range("A33").Value = "오늘은 부모님께 전화 한 통 드리는 것이 어떨까요"
$ws.Range("A34").Value = "한 번 뿐인 인생인데 무엇을 망설이나요"

# Restore view state: 90% zoom, selection on A20, default top-left cell
$ws.Range("A20").Select()
$excel.ActiveWindow.Zoom = 90

